{"js": "// Update the date line and each division problem in the practice table.\n// Old -> New text pairs (every value is unique in the document, so a\n// matchCase/matchWholeWord search safely targets a single run each).\nconst replacements = [\n  [\"2024-03-27 Wednesday\", \"2024-03-28 Thursday\"],\n  [\"521\u00f73=\", \"410\u00f72=\"],\n  [\"669\u00f76=\", \"576\u00f72=\"],\n  [\"708\u00f77=\", \"735\u00f73=\"],\n  [\"967\u00f74=\", \"344\u00f72=\"],\n  [\"255\u00f78=\", \"390\u00f79=\"],\n  [\"788\u00f79=\", \"833\u00f75=\"],\n  [\"853\u00f73=\", \"581\u00f79=\"],\n  [\"110\u00f72=\", \"113\u00f77=\"],\n  [\"465\u00f79=\", \"668\u00f73=\"],\n  [\"573\u00f77=\", \"445\u00f76=\"],\n  [\"755\u00f78=\", \"779\u00f79=\"],\n  [\"117\u00f73=\", \"924\u00f78=\"],\n  [\"581\u00f72=\", \"492\u00f77=\"],\n  [\"198\u00f75=\", \"976\u00f79=\"],\n  [\"101\u00f76=\", \"867\u00f72=\"],\n  [\"497\u00f76=\", \"668\u00f73=\"],\n  [\"494\u00f78=\", \"934\u00f74=\"],\n  [\"167\u00f74=\", \"313\u00f77=\"],\n  [\"305\u00f77=\", \"243\u00f78=\"],\n  [\"997\u00f76=\", \"627\u00f72=\"],\n  [\"645\u00f78=\", \"306\u00f74=\"],\n  [\"214\u00f75=\", \"746\u00f79=\"],\n  [\"249\u00f78=\", \"202\u00f72=\"],\n  [\"194\u00f72=\", \"209\u00f75=\"],\n  [\"190\u00f76=\", \"951\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each division problem in the practice table.\n# Old -> New text pairs (every value is unique in the document, so a\n# case-sensitive whole-document Find/Replace safely targets a single run\n# each).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-27 Wednesday\", \"2024-03-28 Thursday\"),\n    @(\"521\u00f73=\", \"410\u00f72=\"),\n    @(\"669\u00f76=\", \"576\u00f72=\"),\n    @(\"708\u00f77=\", \"735\u00f73=\"),\n    @(\"967\u00f74=\", \"344\u00f72=\"),\n    @(\"255\u00f78=\", \"390\u00f79=\"),\n    @(\"788\u00f79=\", \"833\u00f75=\"),\n    @(\"853\u00f73=\", \"581\u00f79=\"),\n    @(\"110\u00f72=\", \"113\u00f77=\"),\n    @(\"465\u00f79=\", \"668\u00f73=\"),\n    @(\"573\u00f77=\", \"445\u00f76=\"),\n    @(\"755\u00f78=\", \"779\u00f79=\"),\n    @(\"117\u00f73=\", \"924\u00f78=\"),\n    @(\"581\u00f72=\", \"492\u00f77=\"),\n    @(\"198\u00f75=\", \"976\u00f79=\"),\n    @(\"101\u00f76=\", \"867\u00f72=\"),\n    @(\"497\u00f76=\", \"668\u00f73=\"),\n    @(\"494\u00f78=\", \"934\u00f74=\"),\n    @(\"167\u00f74=\", \"313\u00f77=\"),\n    @(\"305\u00f77=\", \"243\u00f78=\"),\n    @(\"997\u00f76=\", \"627\u00f72=\"),\n    @(\"645\u00f78=\", \"306\u00f74=\"),\n    @(\"214\u00f75=\", \"746\u00f79=\"),\n    @(\"249\u00f78=\", \"202\u00f72=\"),\n    @(\"194\u00f72=\", \"209\u00f75=\"),\n    @(\"190\u00f76=\", \"951\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
